$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all text-bearing columns so Excel does not
# auto-convert numeric-looking or date-looking strings.
$ws.Range("A2:C18").NumberFormat = "@"
$ws.Range("E2:E18").NumberFormat = "@"
$ws.Range("H2:H18").NumberFormat = "@"
$ws.Range("J2:J18").NumberFormat = "@"
$ws.Range("M2:N18").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '6262'
$ws.Range("B2").Value = '6/25/2025'
$ws.Range("C2").Value = 'MIGUELETES 1330'
$ws.Range("E2").Value = 'ICD30465943'
$ws.Range("F2").Value = 'Optical Power'
$ws.Range("G2").Value = 'Pendiente'
$ws.Range("H2").Value = 'Cables en panza'
$ws.Range("J2").Value = '{"direccionesNormalizadas": [{"altura": 1330, "cod_calle": 13079, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.440291", "y": "-34.562841"}, "direccion": "MIGUELETES 1330, CABA", "nombre_calle": "MIGUELETES", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M2").Value = 'Colegiales'
$ws.Range("N2").Value = 'Capital Norte'
$ws.Range("D2").Value = 14
$ws.Range("I2").Value = 1
$ws.Range("K2").Value = -58.440291
$ws.Range("L2").Value = -34.562841

# Row 3
$ws.Range("A3").Value = '4756 '
$ws.Range("B3").Value = '12/11/2025'
$ws.Range("C3").Value = 'GARAY, JUAN DE AV. 799'
$ws.Range("E3").Value = '01831840 '
$ws.Range("F3").Value = 'Optical Power'
$ws.Range("G3").Value = 'Pendiente'
$ws.Range("H3").Value = 'tendido bajo'
$ws.Range("J3").Value = '{"direccionesNormalizadas": [{"altura": 799, "cod_calle": 7026, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.376455", "y": "-34.624886"}, "direccion": "GARAY, JUAN DE AV. 799, CABA", "nombre_calle": "GARAY, JUAN DE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M3").Value = 'San Telmo'
$ws.Range("N3").Value = 'Capital Sur'
$ws.Range("D3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("K3").Value = -58.376455
$ws.Range("L3").Value = -34.624886

# Row 4
$ws.Range("A4").Value = '8016'
$ws.Range("B4").Value = '1/8/2026'
$ws.Range("C4").Value = 'RIVADAVIA AV. 2560'
$ws.Range("E4").Value = 'Pendiente ADM'
$ws.Range("F4").Value = 'Optical Power'
$ws.Range("G4").Value = 'Pendiente'
$ws.Range("H4").Value = 'cables sueltos'
$ws.Range("J4").Value = '{"direccionesNormalizadas": [{"altura": 2560, "cod_calle": 19046, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.402591", "y": "-34.610010"}, "direccion": "RIVADAVIA AV. 2560, CABA", "nombre_calle": "RIVADAVIA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M4").Value = 'Almagro'
$ws.Range("N4").Value = 'Capital Sur'
$ws.Range("D4").Value = 3
$ws.Range("I4").Value = 1
$ws.Range("K4").Value = -58.402591
$ws.Range("L4").Value = -34.61001

# Row 5
$ws.Range("A5").Value = '8004'
$ws.Range("B5").Value = '1/8/2026'
$ws.Range("C5").Value = 'LA PAMPA 3650'
$ws.Range("E5").Value = 'Pendiente ADM'
$ws.Range("F5").Value = 'Optical Power'
$ws.Range("G5").Value = 'Pendiente'
$ws.Range("H5").Value = 'caja sobre rejas'
$ws.Range("J5").Value = '{"direccionesNormalizadas": [{"altura": 3650, "cod_calle": 12168, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.466489", "y": "-34.572064"}, "direccion": "LA PAMPA 3650, CABA", "nombre_calle": "LA PAMPA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M5").Value = 'Colegiales'
$ws.Range("N5").Value = 'Capital Norte'
$ws.Range("D5").Value = 13
$ws.Range("I5").Value = 1
$ws.Range("K5").Value = -58.466489
$ws.Range("L5").Value = -34.572064

# Row 6
$ws.Range("A6").Value = '8030'
$ws.Range("B6").Value = '1/8/2026'
$ws.Range("C6").Value = 'SAN JUAN AV. 4267'
$ws.Range("E6").Value = 'Pendiente ADM'
$ws.Range("F6").Value = 'Optical Power'
$ws.Range("G6").Value = 'Pendiente'
$ws.Range("H6").Value = 'cables en panza'
$ws.Range("J6").Value = '{"direccionesNormalizadas": [{"altura": 4267, "cod_calle": 20040, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.425002", "y": "-34.626654"}, "direccion": "SAN JUAN AV. 4267, CABA", "nombre_calle": "SAN JUAN AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M6").Value = 'Boedo'
$ws.Range("N6").Value = 'Capital Sur'
$ws.Range("D6").Value = 5
$ws.Range("I6").Value = 1
$ws.Range("K6").Value = -58.425002
$ws.Range("L6").Value = -34.626654

# Row 7
$ws.Range("A7").Value = '8029'
$ws.Range("B7").Value = '1/8/2026'
$ws.Range("C7").Value = 'MALVINAS ARGENTINAS 49'
$ws.Range("E7").Value = 'Pendiente ADM'
$ws.Range("F7").Value = 'Optical Power'
$ws.Range("G7").Value = 'Pendiente'
$ws.Range("H7").Value = 'cable en panza, cortado'
$ws.Range("J7").Value = '{"direccionesNormalizadas": [{"altura": 49, "cod_calle": 13020, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.451796", "y": "-34.625377"}, "direccion": "MALVINAS ARGENTINAS 49, CABA", "nombre_calle": "MALVINAS ARGENTINAS", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M7").Value = 'Boedo'
$ws.Range("N7").Value = 'Capital Sur'
$ws.Range("D7").Value = 6
$ws.Range("I7").Value = 1
$ws.Range("K7").Value = -58.451796
$ws.Range("L7").Value = -34.625377

# Row 8
$ws.Range("A8").Value = '8048'
$ws.Range("B8").Value = '1/8/2026'
$ws.Range("C8").Value = 'ARANGUREN, JUAN F., DR. 4355'
$ws.Range("E8").Value = 'Pendiente ADM'
$ws.Range("F8").Value = 'Optical Power'
$ws.Range("G8").Value = 'Pendiente'
$ws.Range("H8").Value = 'riesgo de caida de equipo'
$ws.Range("J8").Value = '{"direccionesNormalizadas": [{"altura": 4355, "cod_calle": 1094, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.489757", "y": "-34.630158"}, "direccion": "ARANGUREN, JUAN F., DR. 4355, CABA", "nombre_calle": "ARANGUREN, JUAN F., DR.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M8").Value = 'Devoto'
$ws.Range("N8").Value = 'Capital Norte'
$ws.Range("D8").Value = 10
$ws.Range("I8").Value = 1
$ws.Range("K8").Value = -58.489757
$ws.Range("L8").Value = -34.630158

# Row 9
$ws.Range("A9").Value = 'Z7'
$ws.Range("B9").Value = '1/12/2026'
$ws.Range("C9").Value = 'JULIAN ALVAREZ 928'
$ws.Range("E9").Value = 'Pendiente ADM'
$ws.Range("F9").Value = 'Optical Power'
$ws.Range("G9").Value = 'Pendiente'
$ws.Range("H9").Value = 'cables colgando'
$ws.Range("J9").Value = '{"direccionesNormalizadas": [{"altura": 928, "cod_calle": 1057, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.430574", "y": "-34.597043"}, "direccion": "ALVAREZ, JULIAN 928, CABA", "nombre_calle": "ALVAREZ, JULIAN", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M9").Value = 'Palermo'
$ws.Range("N9").Value = 'Capital Sur'
$ws.Range("D9").Value = 15
$ws.Range("I9").Value = 1
$ws.Range("K9").Value = -58.430574
$ws.Range("L9").Value = -34.597043

# Row 10
$ws.Range("A10").Value = 'Z1'
$ws.Range("B10").Value = '1/12/2026'
$ws.Range("C10").Value = 'LAVALLEJA 990'
$ws.Range("E10").Value = 'Pendiente ADM'
$ws.Range("F10").Value = 'Optical Power'
$ws.Range("G10").Value = 'Pendiente'
$ws.Range("H10").Value = 'cable cortado'
$ws.Range("J10").Value = '{"direccionesNormalizadas": [{"altura": 990, "cod_calle": 12090, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.429174", "y": "-34.597484"}, "direccion": "LAVALLEJA 990, CABA", "nombre_calle": "LAVALLEJA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M10").Value = 'Palermo'
$ws.Range("N10").Value = 'Capital Sur'
$ws.Range("D10").Value = 15
$ws.Range("I10").Value = 1
$ws.Range("K10").Value = -58.429174
$ws.Range("L10").Value = -34.597484

# Row 11
$ws.Range("A11").Value = 'Z2'
$ws.Range("B11").Value = '1/12/2026'
$ws.Range("C11").Value = 'LAVALLEJA 1030'
$ws.Range("E11").Value = 'Pendiente ADM'
$ws.Range("F11").Value = 'Optical Power'
$ws.Range("G11").Value = 'Pendiente'
$ws.Range("H11").Value = 'altura insufciente'
$ws.Range("J11").Value = '{"direccionesNormalizadas": [{"altura": 1030, "cod_calle": 12090, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.428659", "y": "-34.597153"}, "direccion": "LAVALLEJA 1030, CABA", "nombre_calle": "LAVALLEJA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M11").Value = 'Palermo'
$ws.Range("N11").Value = 'Capital Sur'
$ws.Range("D11").Value = 15
$ws.Range("I11").Value = 1
$ws.Range("K11").Value = -58.428659
$ws.Range("L11").Value = -34.597153

# Row 12
$ws.Range("A12").Value = 'Z3'
$ws.Range("B12").Value = '1/12/2026'
$ws.Range("C12").Value = 'QUILMES 181'
$ws.Range("E12").Value = 'Pendiente ADM'
$ws.Range("F12").Value = 'Optical Power'
$ws.Range("G12").Value = 'Pendiente'
$ws.Range("H12").Value = 'cable colgando'
$ws.Range("J12").Value = '{"direccionesNormalizadas": [{"altura": 181, "cod_calle": 18007, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.418079", "y": "-34.641658"}, "direccion": "QUILMES 181, CABA", "nombre_calle": "QUILMES", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M12").Value = 'Boedo'
$ws.Range("N12").Value = 'Capital Sur'
$ws.Range("D12").Value = 4
$ws.Range("I12").Value = 1
$ws.Range("K12").Value = -58.418079
$ws.Range("L12").Value = -34.641658

# Row 13
$ws.Range("A13").Value = 'Z4'
$ws.Range("B13").Value = '1/12/2026'
$ws.Range("C13").Value = 'QUILMES 256'
$ws.Range("E13").Value = 'Pendiente ADM'
$ws.Range("F13").Value = 'Optical Power'
$ws.Range("G13").Value = 'Pendiente'
$ws.Range("H13").Value = 'cables colgando'
$ws.Range("J13").Value = '{"direccionesNormalizadas": [{"altura": 256, "cod_calle": 18007, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.417736", "y": "-34.642492"}, "direccion": "QUILMES 256, CABA", "nombre_calle": "QUILMES", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M13").Value = 'Boedo'
$ws.Range("N13").Value = 'Capital Sur'
$ws.Range("D13").Value = 4
$ws.Range("I13").Value = 1
$ws.Range("K13").Value = -58.417736
$ws.Range("L13").Value = -34.642492

# Row 14
$ws.Range("A14").Value = 'Z5'
$ws.Range("B14").Value = '1/12/2026'
$ws.Range("C14").Value = 'SERRANO 1074'
$ws.Range("E14").Value = 'Pendiente ADM'
$ws.Range("F14").Value = 'Optical Power'
$ws.Range("G14").Value = 'Pendiente'
$ws.Range("H14").Value = 'cable colgando'
$ws.Range("J14").Value = '{"direccionesNormalizadas": [{"altura": 1074, "cod_calle": 20090, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.435899", "y": "-34.592365"}, "direccion": "SERRANO 1074, CABA", "nombre_calle": "SERRANO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M14").Value = 'Palermo'
$ws.Range("N14").Value = 'Capital Sur'
$ws.Range("D14").Value = 15
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = -58.435899
$ws.Range("L14").Value = -34.592365

# Row 15
$ws.Range("A15").Value = 'Z6'
$ws.Range("B15").Value = '1/12/2026'
$ws.Range("C15").Value = 'CORRIENTES AV. 5733'
$ws.Range("E15").Value = 'Pendiente ADM'
$ws.Range("F15").Value = 'Optical Power'
$ws.Range("G15").Value = 'Pendiente'
$ws.Range("H15").Value = 'rg11 colgando'
$ws.Range("J15").Value = '{"direccionesNormalizadas": [{"altura": 5733, "cod_calle": 3174, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.443697", "y": "-34.595202"}, "direccion": "CORRIENTES AV. 5733, CABA", "nombre_calle": "CORRIENTES AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M15").Value = 'Palermo'
$ws.Range("N15").Value = 'Capital Sur'
$ws.Range("D15").Value = 15
$ws.Range("I15").Value = 1
$ws.Range("K15").Value = -58.443697
$ws.Range("L15").Value = -34.595202

# Row 16
$ws.Range("A16").Value = 'Z8'
$ws.Range("B16").Value = '1/12/2026'
$ws.Range("C16").Value = 'ALVAREZ JONTE AV. 1808'
$ws.Range("E16").Value = 'Pendiente ADM'
$ws.Range("F16").Value = 'Optical Power'
$ws.Range("G16").Value = 'Pendiente'
$ws.Range("H16").Value = 'cdo colgando'
$ws.Range("J16").Value = '{"direccionesNormalizadas": [{"altura": 1808, "cod_calle": 1056, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.469350", "y": "-34.602268"}, "direccion": "ALVAREZ JONTE 1808, CABA", "nombre_calle": "ALVAREZ JONTE", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M16").Value = 'Paternal'
$ws.Range("N16").Value = 'Capital Norte'
$ws.Range("D16").Value = 11
$ws.Range("I16").Value = 1
$ws.Range("K16").Value = -58.46935
$ws.Range("L16").Value = -34.602268

# Row 17
$ws.Range("A17").Value = 'S01268691'
$ws.Range("B17").Value = '1/12/2026'
$ws.Range("C17").Value = 'ESTADOS UNIDOS 1943'
$ws.Range("E17").Value = 'Pendiente ADM'
$ws.Range("F17").Value = 'Optical Power'
$ws.Range("G17").Value = 'Pendiente'
$ws.Range("H17").Value = 'cables colgando'
$ws.Range("J17").Value = '{"direccionesNormalizadas": [{"altura": 1943, "cod_calle": 5087, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.393711", "y": "-34.619471"}, "direccion": "ESTADOS UNIDOS 1943, CABA", "nombre_calle": "ESTADOS UNIDOS", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M17").Value = 'San Telmo'
$ws.Range("N17").Value = 'Capital Sur'
$ws.Range("D17").Value = 3
$ws.Range("I17").Value = 1
$ws.Range("K17").Value = -58.393711
$ws.Range("L17").Value = -34.619471

# Row 18
$ws.Range("A18").Value = 'S01335742'
$ws.Range("B18").Value = '1/12/2026'
$ws.Range("C18").Value = 'Tinogasta 5182'
$ws.Range("E18").Value = 'Pendiente ADM'
$ws.Range("F18").Value = 'Optical Power'
$ws.Range("G18").Value = 'Pendiente'
$ws.Range("H18").Value = 'cables bajos'
$ws.Range("J18").Value = '{"direccionesNormalizadas": [{"altura": 5182, "cod_calle": 21032, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.519521", "y": "-34.615739"}, "direccion": "TINOGASTA 5182, CABA", "nombre_calle": "TINOGASTA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("M18").Value = 'Devoto'
$ws.Range("N18").Value = 'Capital Norte'
$ws.Range("D18").Value = 11
$ws.Range("I18").Value = 1
$ws.Range("K18").Value = -58.519521
$ws.Range("L18").Value = -34.615739
